$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 108, shifting existing rows 108-118 down to 109-119
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with data
$ws.Cells.Item(108, 1).Value = 9
$ws.Cells.Item(108, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(108, 3).Value = "Metropolitana"
$ws.Cells.Item(108, 4).Value = 44769
$ws.Cells.Item(108, 5).Value = 13
$ws.Cells.Item(108, 6).Value = 100112022
$ws.Cells.Item(108, 7).Value = "Arveja Verde"
$ws.Cells.Item(108, 8).Value = "Perfection"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 25
$ws.Cells.Item(108, 11).Value = 43000
$ws.Cells.Item(108, 12).Value = 43000
$ws.Cells.Item(108, 13).Value = 43000
$ws.Cells.Item(108, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(108, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(108, 16).Value = 1720
$ws.Cells.Item(108, 17).Value = 25
$ws.Cells.Item(108, 18).Value = "Hortaliza"
